# Add a new "Hello world (Wingdings smiley)" line, authored with the
# Croatian (hr-HR) proofing language, ahead of the existing (empty)
# paragraph that holds the _GoBack bookmark.

$d = $word.ActiveDocument

# The document currently has a single, empty paragraph. Stamp the
# paragraph mark's language first (hr-HR) -- with no runs yet present,
# setting LanguageID on the whole paragraph range lands on the
# paragraph-mark run properties (w:pPr/w:rPr/w:lang), exactly as Word
# does when you switch the language for an empty paragraph before
# typing into it.
$p = $d.Paragraphs(1)
$p.Range.LanguageID = "hr-HR"

# Build the new content as literal OOXML: "Hello" and "world" are
# wrapped in spell-check proofErr markers (they are not Croatian
# words), separated by plain spaces, followed by a Wingdings smiley
# symbol run -- all tagged with the hr-HR run language, matching what
# Word inserts when proofing a typed sentence.
$xml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:proofErr w:type='spellStart'/>
  <w:r>
    <w:rPr>
      <w:lang w:val='hr-HR'/>
    </w:rPr>
    <w:t>Hello</w:t>
  </w:r>
  <w:proofErr w:type='spellEnd'/>
  <w:r>
    <w:rPr>
      <w:lang w:val='hr-HR'/>
    </w:rPr>
    <w:t xml:space='preserve'> </w:t>
  </w:r>
  <w:proofErr w:type='spellStart'/>
  <w:r>
    <w:rPr>
      <w:lang w:val='hr-HR'/>
    </w:rPr>
    <w:t>world</w:t>
  </w:r>
  <w:proofErr w:type='spellEnd'/>
  <w:r>
    <w:rPr>
      <w:lang w:val='hr-HR'/>
    </w:rPr>
    <w:t xml:space='preserve'> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val='hr-HR'/>
    </w:rPr>
    <w:sym w:font='Wingdings' w:char='F04A'/>
  </w:r>
</w:p>
"@

# InsertXML operates at paragraph granularity: it inserts the supplied
# markup as a brand-new paragraph just before the (still empty)
# original paragraph, rather than splicing runs into it.
$r = $d.Range(0, 0)
$r.InsertXML($xml)

# Merge the newly-inserted paragraph with the original one (which
# still carries the pPr/lang we stamped, and the _GoBack bookmark) by
# deleting the paragraph mark between them -- same effect as pressing
# Delete at the end of the first line to join it with the next.
$p1 = $d.Paragraphs(1)
$endOfP1 = $p1.Range.End
$d.Range($endOfP1 - 1, $endOfP1).Delete()
